$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.837.71"
$ws.Range("E2").Value = "  +8.36%  "

$ws.Range("D3").Value = "3.216.05"
$ws.Range("E3").Value = "  +3.63%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "397.07"
$ws.Range("E5").Value = "  +2.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.63"
$ws.Range("E6").Value = "  +5.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.555"
$ws.Range("E7").Value = "  +2.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  +6.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.29"
$ws.Range("E10").Value = "  +6.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  +6.41%  "

$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("D13").Value = "3.720.73"
$ws.Range("E13").Value = "  +3.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.09"
$ws.Range("E14").Value = "  +2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.08"
$ws.Range("E15").Value = "  +3.38%  "

$ws.Range("E16").Value = "  +6.15%  "

$ws.Range("D17").Value = "3.193.90"
$ws.Range("E17").Value = "  +3.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.67"
$ws.Range("E18").Value = "  -2.37%  "

$ws.Range("D19").Value = "55.695.05"
$ws.Range("E19").Value = "  +7.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.37"
$ws.Range("E20").Value = "  +3.89%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000102"
$ws.Range("E21").Value = "  +6.10%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.10"
$ws.Range("E22").Value = "  +5.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "304.57"
$ws.Range("E23").Value = "  +14.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.32"
$ws.Range("E24").Value = "  +7.81%  "

$ws.Range("E25").Value = "  +1.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.25"
$ws.Range("E26").Value = "  +1.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.28"
$ws.Range("E27").Value = "  +4.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.53"
$ws.Range("E28").Value = "  +4.55%  "

$ws.Range("E29").Value = "  +4.69%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.111"
$ws.Range("E31").Value = "  +4.31%  "

$ws.Range("E32").Value = "  +9.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0495"
$ws.Range("E33").Value = "  +3.25%  "

$ws.Range("E34").Value = "  +2.63%  "

$ws.Range("E35").Value = "  +2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.18"
$ws.Range("E36").Value = "  +2.30%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.10"
$ws.Range("E37").Value = "  +23.82%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  +4.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "134.82"
$ws.Range("E40").Value = "  +4.46%  "

$ws.Range("E42").Value = "  +2.47%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.16"
$ws.Range("E43").Value = "  +3.45%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.287"
$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.119"
$ws.Range("E45").Value = "  +2.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.23"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.14"
$ws.Range("E47").Value = "  +46.05%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  -0.99%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.09"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("D50").Value = "2.138.60"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0362"
$ws.Range("E51").Value = "  +9.92%  "
